# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Lane Late - Primera / Segunda, dated 44508)
# at rows 411-412, pushing the existing rows 411-523 down to 413-525.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 411 so every following row
# (formerly 411..523) shifts down by two, ending at 413..525.
$ws.Rows("411:412").Insert()

# --- New row 411: Lane Late / Primera ---
$ws.Cells.Item(411, 1).Value = 3
$ws.Cells.Item(411, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(411, 3).Value = "Coquimbo"
$ws.Cells.Item(411, 4).Value = 44508
$ws.Cells.Item(411, 5).Value = 5
$ws.Cells.Item(411, 6).Value = "Fruta"
$ws.Cells.Item(411, 7).Value = 100102
$ws.Cells.Item(411, 8).Value = "Cítricos"
$ws.Cells.Item(411, 9).Value = 100102005
$ws.Cells.Item(411, 10).Value = "Naranja"
$ws.Cells.Item(411, 11).Value = "Lane Late"
$ws.Cells.Item(411, 12).Value = "Primera"
$ws.Cells.Item(411, 13).Value = 188
$ws.Cells.Item(411, 14).Value = 5000
$ws.Cells.Item(411, 15).Value = 6000
$ws.Cells.Item(411, 16).Value = 5479
$ws.Cells.Item(411, 17).Value = "`$/malla 13 kilos"
$ws.Cells.Item(411, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(411, 19).Value = 421
$ws.Cells.Item(411, 20).Value = 13

# --- New row 412: Lane Late / Segunda ---
$ws.Cells.Item(412, 1).Value = 3
$ws.Cells.Item(412, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(412, 3).Value = "Coquimbo"
$ws.Cells.Item(412, 4).Value = 44508
$ws.Cells.Item(412, 5).Value = 5
$ws.Cells.Item(412, 6).Value = "Fruta"
$ws.Cells.Item(412, 7).Value = 100102
$ws.Cells.Item(412, 8).Value = "Cítricos"
$ws.Cells.Item(412, 9).Value = 100102005
$ws.Cells.Item(412, 10).Value = "Naranja"
$ws.Cells.Item(412, 11).Value = "Lane Late"
$ws.Cells.Item(412, 12).Value = "Segunda"
$ws.Cells.Item(412, 13).Value = 174
$ws.Cells.Item(412, 14).Value = 4000
$ws.Cells.Item(412, 15).Value = 4500
$ws.Cells.Item(412, 16).Value = 4250
$ws.Cells.Item(412, 17).Value = "`$/malla 13 kilos"
$ws.Cells.Item(412, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(412, 19).Value = 327
$ws.Cells.Item(412, 20).Value = 13
